{"js": "// The \"Host: AWS\" paragraph currently reads (in part):\n//   \"... Our code will also clean, normalize, and standardize data prior to modeling.\"\n// It needs to become:\n//   \"... Our code will also clean our data prior to modeling.\"\n// and the part that was trimmed out (\"normalize, and standardize data prior\n// to modeling\") reappears as a brand-new leading sentence on the very next\n// content paragraph (\"Model: Logistical Regression\" bullet body), turning\n// \"Initialize, train, and evaluate our model. ...\" into\n// \"Normalize and standardize data prior to modeling. Initialize, train, and\n// evaluate our model. ...\".\n\nconst body = context.document.body;\n\n// 1) Locate the exact run of text that needs trimming down and rewrite it.\nconst oldSentence = body.search(\n  \"lean, normalize, and standardize data prior to modeling\",\n  { matchCase: true }\n);\noldSentence.load(\"text\");\nawait context.sync();\n\nif (oldSentence.items.length > 0) {\n  oldSentence.items[0].insertText(\"lean our data prior to modeling\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Locate the paragraph that starts with \"Initialize, train, and evaluate\"\n// and prepend the new \"Normalize and standardize data prior to modeling. \"\n// sentence in front of it.\nconst nextParaStart = body.search(\"Initialize, train, and evaluate\", { matchCase: true });\nnextParaStart.load(\"text\");\nawait context.sync();\n\nif (nextParaStart.items.length > 0) {\n  nextParaStart.items[0].insertText(\n    \"Normalize and standardize data prior to modeling. \",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n}\n", "ps1": "# The \"Host: AWS\" paragraph currently reads (in part):\n#   \"... Our code will also clean, normalize, and standardize data prior to modeling.\"\n# It needs to become:\n#   \"... Our code will also clean our data prior to modeling.\"\n# and the trimmed-out part (\"normalize, and standardize data prior to\n# modeling\") reappears as a brand-new leading sentence on the very next\n# content paragraph (\"Model: Logistical Regression\" bullet body), turning\n# \"Initialize, train, and evaluate our model. ...\" into\n# \"Normalize and standardize data prior to modeling. Initialize, train, and\n# evaluate our model. ...\".\n\n$d = $word.ActiveDocument\n\n# 1) Trim the \"Host: AWS\" sentence down to \"... clean our data prior to modeling.\"\n$range1 = $d.Content\n$found1 = $range1.Find.Execute(\n    \"lean, normalize, and standardize data prior to modeling\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"lean our data prior to modeling\",\n    1\n)\nWrite-Output \"replaced sentence trim: $found1\"\n\n# 2) Prepend the new sentence to the \"Model: Logistical Regression\" paragraph.\n$range2 = $d.Content\n$found2 = $range2.Find.Execute(\n    \"Initialize, train, and evaluate\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Normalize and standardize data prior to modeling. Initialize, train, and evaluate\",\n    1\n)\nWrite-Output \"prepended new sentence: $found2\"\n"}
